# Updates cryptos list figures (price / 1h volume change) and reorders
# a couple of coin rows, per the scraper's latest run.
#
# Note: Price values in column D are stored as plain TEXT in the workbook
# (e.g. "62.170.31"), not numbers. Values that look like a genuine number
# (e.g. "555.88") would otherwise be auto-converted to a numeric value by
# Excel, so those are written with a leading apostrophe to force them to
# stay text, matching the original cell typing. Values that already
# contain two dots (e.g. "62.170.31") are never interpreted as numbers,
# so no apostrophe is needed for those.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple price / percentage updates -------------------------------------------------
$ws.Range("D2").Value  = "62.170.31"
$ws.Range("E2").Value  = "  +2.08%  "

$ws.Range("D3").Value  = "2.419.99"
$ws.Range("E3").Value  = "  +2.78%  "

$ws.Range("D5").Value  = "'555.88"
$ws.Range("E5").Value  = "  +1.88%  "

$ws.Range("D6").Value  = "'142.97"
$ws.Range("E6").Value  = "  +4.21%  "

$ws.Range("E8").Value  = "  +1.42%  "

$ws.Range("D9").Value  = "2.420.78"
$ws.Range("E9").Value  = "  +2.85%  "

$ws.Range("E10").Value = "  +3.80%  "

$ws.Range("E12").Value = "  +1.47%  "

$ws.Range("E13").Value = "  +1.23%  "

$ws.Range("D14").Value = "'26.19"
$ws.Range("E14").Value = "  +5.88%  "

$ws.Range("E15").Value = "  +7.40%  "

$ws.Range("D16").Value = "2.861.14"
$ws.Range("E16").Value = "  +3.01%  "

$ws.Range("D17").Value = "62.167.36"
$ws.Range("E17").Value = "  +1.82%  "

$ws.Range("D18").Value = "2.421.33"
$ws.Range("E18").Value = "  +3.06%  "

$ws.Range("D19").Value = "'11.07"
$ws.Range("E19").Value = "  +3.86%  "

$ws.Range("D20").Value = "'4.19"
$ws.Range("E20").Value = "  +1.59%  "

$ws.Range("D21").Value = "'324.29"
$ws.Range("E21").Value = "  +1.09%  "

$ws.Range("E22").Value = "  +2.04%  "

$ws.Range("E23").Value = "  +0.20%  "

$ws.Range("D24").Value = "'64.94"
$ws.Range("E24").Value = "  +2.37%  "

$ws.Range("E25").Value = "  +4.58%  "

$ws.Range("D26").Value = "'9.11"
$ws.Range("E26").Value = "  +7.45%  "

$ws.Range("D27").Value = "'578.59"
$ws.Range("E27").Value = "  +16.43%  "

$ws.Range("D28").Value = "2.542.39"
$ws.Range("E28").Value = "  +3.05%  "

$ws.Range("E29").Value = "  -0.05%  "

$ws.Range("D30").Value = "'8.38"
$ws.Range("E30").Value = "  +4.17%  "

$ws.Range("E31").Value = "  +7.06%  "

$ws.Range("E32").Value = "  +5.27%  "

$ws.Range("E33").Value = "  +1.29%  "

$ws.Range("E34").Value = "  +3.54%  "

$ws.Range("E35").Value = "  +2.70%  "

$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  +0.14%  "

# --- Row 37/38 swap: NEARProtocol <-> RenderToken ---------------------------------------
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D37").Value = "'5.67"
$ws.Range("E37").Value = "  +7.80%  "

$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").Value = "'4.82"
$ws.Range("E38").Value = "  +3.53%  "

$ws.Range("E39").Value = "  +1.69%  "

# --- Row 40/41 swap: Stacks <-> EthereumClassic ------------------------------------------
$ws.Range("B40").Value = "EthereumClassic"
$ws.Range("C40").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D40").Value = "'18.75"
$ws.Range("E40").Value = "  +1.08%  "

$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'1.87"
$ws.Range("E41").Value = "  +2.30%  "

# --- Remaining simple updates -------------------------------------------------------------
$ws.Range("D42").Value = "'148.01"
$ws.Range("E42").Value = "  +2.22%  "

$ws.Range("E43").Value = "  +0.06%  "

$ws.Range("D44").Value = "'41.70"
$ws.Range("E44").Value = "  +2.73%  "

$ws.Range("E45").Value = "  +11.28%  "

$ws.Range("D46").Value = "'150.88"
$ws.Range("E46").Value = "  +5.76%  "

$ws.Range("E47").Value = "  +1.31%  "

$ws.Range("D48").Value = "'0.0544"
$ws.Range("E48").Value = "  +5.31%  "

$ws.Range("D49").Value = "'20.34"
$ws.Range("E49").Value = "  +5.92%  "

$ws.Range("E50").Value = "  +3.42%  "

$ws.Range("D51").Value = "'0.0918"
$ws.Range("E51").Value = "  +1.62%  "
